$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")
$cos = $ws.ChartObjects()
$co = $cos.Item(1)
$chart = $co.Chart
$chart.HasTitle = $true
$ct = $chart.ChartTitle
$ct.Text = "2 seconds"
Write-Output $ct.Top
Write-Output $ct.Left
